$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated forecast values (rows 2-7) re-generated from R ("Aggiorno file need_to_buy.xlsx da R")
$ws.Range("B2").Value = 7920.20922157451
$ws.Range("C2").Value = 7061.02150765286
$ws.Range("D2").Value = 10924
$ws.Range("E2").Value = 3483.83992799748
$ws.Range("F2").Value = -15.7974401812359

$ws.Range("B3").Value = 8101.48141026903
$ws.Range("C3").Value = 7821.02291501942
$ws.Range("E3").Value = 3799.02418697549
$ws.Range("F3").Value = 109.001962583121

$ws.Range("B4").Value = 8482.76857383341
$ws.Range("C4").Value = 8228.93224265405
$ws.Range("E4").Value = 4624.02131548604
$ws.Range("F4").Value = 160.373064922504

$ws.Range("B5").Value = 8428.76975398309
$ws.Range("C5").Value = 8372.54875692511
$ws.Range("E5").Value = 4639.1762502409
$ws.Range("F5").Value = 166.98854196525

$ws.Range("B6").Value = 7969.35648312576
$ws.Range("C6").Value = 7457.0926512824
$ws.Range("E6").Value = 4281.95080701943
$ws.Range("F6").Value = 113.96014409591

$ws.Range("B7").Value = 6865.32977546962
$ws.Range("C7").Value = 7239.29943752965
$ws.Range("E7").Value = 4252.91975834716
$ws.Range("F7").Value = 103.6757998282
